$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Locate the two paragraphs that make up the "Capture our own dataset..."
#    / "Train on our dataset..." checklist items (by their text, so the
#    script does not depend on fragile paragraph indices).
# ---------------------------------------------------------------------------
$captureIdx = -1
$trainIdx = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $t = $p.Range.Text
    if ($t -like "*Capture our own dataset of images from*") {
        $captureIdx = $idx
    }
    if ($t -like "*Train on our dataset to see some preliminary results*") {
        $trainIdx = $idx
    }
}

if ($captureIdx -eq -1 -or $trainIdx -eq -1) {
    throw "Could not locate the target checklist paragraphs (capture=$captureIdx, train=$trainIdx)"
}

$pCapture = $d.Paragraphs.Item($captureIdx)
$pTrain = $d.Paragraphs.Item($trainIdx)

$rangeStart = $pCapture.Range.Start
$rangeEnd = $pTrain.Range.End
$target = $d.Range($rangeStart, $rangeEnd)

# Replacement XML: a new checklist item ("Access NUbots synthetic dataset and
# train on rectilinear images", reusing the original checkbox id and gaining
# the _GoBack bookmark) followed by the original two paragraphs verbatim
# (now carrying a duplicated checkbox id of 1194885365, plus the sdtEndPr
# elements and losing the bookmark that moved to the new first paragraph).
$newItemsXml = '<w:p w:rsidR="00103E84" w:rsidRDefault="00103E84" w:rsidP="00103E84" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:left="709" w:hanging="709"/></w:pPr><w:sdt><w:sdtPr><w:id w:val="141161360"/><w14:checkbox><w14:checked w14:val="0"/><w14:checkedState w14:val="00FE" w14:font="Wingdings"/><w14:uncheckedState w14:val="2610" w14:font="MS Gothic"/></w14:checkbox></w:sdtPr><w:sdtContent><w:r><w:rPr><w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:hint="eastAsia"/></w:rPr><w:t>&#9744;</w:t></w:r></w:sdtContent></w:sdt><w:r><w:tab/><w:t>Access NUbots synthetic dataset and train on rectilinear images</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p w:rsidR="00103E84" w:rsidRDefault="00103E84" w:rsidP="00103E84" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:left="709" w:hanging="709"/></w:pPr><w:sdt><w:sdtPr><w:id w:val="1194885365"/><w14:checkbox><w14:checked w14:val="0"/><w14:checkedState w14:val="00FE" w14:font="Wingdings"/><w14:uncheckedState w14:val="2610" w14:font="MS Gothic"/></w14:checkbox></w:sdtPr><w:sdtEndPr/><w:sdtContent><w:r><w:rPr><w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:hint="eastAsia"/></w:rPr><w:t>&#9744;</w:t></w:r></w:sdtContent></w:sdt><w:r><w:tab/><w:t xml:space="preserve">Capture our own dataset of images from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>igus</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p w:rsidR="00103E84" w:rsidRDefault="00103E84" w:rsidP="00103E84" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:left="709" w:hanging="709"/></w:pPr><w:sdt><w:sdtPr><w:id w:val="1155261351"/><w14:checkbox><w14:checked w14:val="0"/><w14:checkedState w14:val="00FE" w14:font="Wingdings"/><w14:uncheckedState w14:val="2610" w14:font="MS Gothic"/></w14:checkbox></w:sdtPr><w:sdtEndPr/><w:sdtContent><w:r><w:rPr><w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:hint="eastAsia"/></w:rPr><w:t>&#9744;</w:t></w:r></w:sdtContent></w:sdt><w:r><w:tab/><w:t>Train on our dataset to see some preliminary results</w:t></w:r></w:p>'

$target.InsertXML($newItemsXml)

# ---------------------------------------------------------------------------
# 2) Remove the stray <w:lastRenderedPageBreak/> run child further down the
#    document (the paragraph whose entire content is a single tab character).
# ---------------------------------------------------------------------------
$tabParaIdx = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -eq "`t`r") {
        $tabParaIdx = $idx
    }
}

if ($tabParaIdx -eq -1) {
    throw "Could not locate the lone-tab paragraph with the lastRenderedPageBreak run"
}

$pTab = $d.Paragraphs.Item($tabParaIdx)
$tabRange = $pTab.Range
$tabXml = '<w:p w:rsidR="00C33550" w:rsidRDefault="00C33550" w:rsidP="00C33550" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="7710"/></w:tabs></w:pPr><w:r><w:tab/></w:r></w:p>'
$tabRange.InsertXML($tabXml)

Write-Output "Edit complete."
